$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.307.22'
$ws.Range("E2").Value = '  -1.92%  '
$ws.Range("D3").Value = '2.503.27'
$ws.Range("E3").Value = '  -3.91%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.16'
$ws.Range("E5").Value = '  -2.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.49'
$ws.Range("E6").Value = '  -4.69%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").Value = '2.500.26'
$ws.Range("E9").Value = '  -3.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.109'
$ws.Range("E10").Value = '  -7.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.155'
$ws.Range("E11").Value = '  -0.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.44'
$ws.Range("E12").Value = '  -6.29%  '
$ws.Range("E13").Value = '  -5.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.35'
$ws.Range("E14").Value = '  -6.68%  '
$ws.Range("D15").Value = '2.958.72'
$ws.Range("E15").Value = '  -3.85%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '62.194.24'
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000165'
$ws.Range("E17").Value = '  -6.60%  '
$ws.Range("D18").Value = '2.503.10'
$ws.Range("E18").Value = '  -4.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.30'
$ws.Range("E19").Value = '  -5.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.05'
$ws.Range("E20").Value = '  -5.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.25'
$ws.Range("E21").Value = '  -6.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '324.95'
$ws.Range("E22").Value = '  -5.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.50'
$ws.Range("E24").Value = '  -3.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.77'
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000104'
$ws.Range("E26").Value = '  -3.31%  '
$ws.Range("D27").Value = '2.629.75'
$ws.Range("E27").Value = '  -3.85%  '
$ws.Range("E28").Value = '  -3.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.52'
$ws.Range("E29").Value = '  -6.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '540.59'
$ws.Range("E30").Value = '  -7.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.77'
$ws.Range("E32").Value = '  -1.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.152'
$ws.Range("E33").Value = '  -4.92%  '
$ws.Range("E34").Value = '  -6.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").Value = '  -7.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.96'
$ws.Range("E36").Value = '  -8.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.90'
$ws.Range("E37").Value = '  -8.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("E39").Value = '  -5.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.70'
$ws.Range("E40").Value = '  -4.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.27'
$ws.Range("E41").Value = '  -5.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.72'
$ws.Range("E42").Value = '  -7.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.85'
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '149.36'
$ws.Range("E46").Value = '  -4.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.66'
$ws.Range("E47").Value = '  -6.04%  '
$ws.Range("E48").Value = '  -6.48%  '
$ws.Range("E49").Value = '  -7.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.596'
$ws.Range("E50").Value = '  -5.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0956'
$ws.Range("E51").Value = '  -4.41%  '
